$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 223 (pushes existing rows 223-271 down to 225-273),
# matching the weekly data refresh described in the commit message.
$ws.Rows.Item(223).Insert()
$ws.Rows.Item(223).Insert()

# New row 223: Coliflor, Primera, Region Metropolitana, date 44754 (2022-07-12)
$ws.Cells.Item(223, 1).Value = 11
$ws.Cells.Item(223, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(223, 3).Value = "Bíobío"
$ws.Cells.Item(223, 4).Value = 44754
$ws.Cells.Item(223, 5).Value = 8
$ws.Cells.Item(223, 6).Value = 100112008
$ws.Cells.Item(223, 7).Value = "Coliflor"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 2000
$ws.Cells.Item(223, 11).Value = 1100
$ws.Cells.Item(223, 12).Value = 1200
$ws.Cells.Item(223, 13).Value = 1150
$ws.Cells.Item(223, 14).Value = "`$/unidad"
$ws.Cells.Item(223, 15).Value = "Región Metropolitana"
$ws.Cells.Item(223, 16).Value = 1150
$ws.Cells.Item(223, 17).Value = 1
$ws.Cells.Item(223, 18).Value = "Hortaliza"

# New row 224: Coliflor, Segunda, Region Metropolitana, same date 44754
$ws.Cells.Item(224, 1).Value = 11
$ws.Cells.Item(224, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(224, 3).Value = "Bíobío"
$ws.Cells.Item(224, 4).Value = 44754
$ws.Cells.Item(224, 5).Value = 8
$ws.Cells.Item(224, 6).Value = 100112008
$ws.Cells.Item(224, 7).Value = "Coliflor"
$ws.Cells.Item(224, 8).Value = "Sin especificar"
$ws.Cells.Item(224, 9).Value = "Segunda"
$ws.Cells.Item(224, 10).Value = 1000
$ws.Cells.Item(224, 11).Value = 900
$ws.Cells.Item(224, 12).Value = 900
$ws.Cells.Item(224, 13).Value = 900
$ws.Cells.Item(224, 14).Value = "`$/unidad"
$ws.Cells.Item(224, 15).Value = "Región Metropolitana"
$ws.Cells.Item(224, 16).Value = 900
$ws.Cells.Item(224, 17).Value = 1
$ws.Cells.Item(224, 18).Value = "Hortaliza"
